# Insert a new weekly data row for "Arándano (blue)" (Macroferia Regional de Talca)
# at row 25 of Sheet1, pushing the existing rows 25..58 down to 26..59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 25 (this shifts rows 25-58 down to 26-59,
# carrying their formatting with them).
$ws.Rows.Item(25).Insert()

# Populate the newly-inserted row 25 with this week's data.
$ws.Cells.Item(25, 1).Value = 5
$ws.Cells.Item(25, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(25, 3).Value = "Maule"
$ws.Cells.Item(25, 4).Value = 44589
$ws.Cells.Item(25, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(25, 5).Value = 7
$ws.Cells.Item(25, 6).Value = "Fruta"
$ws.Cells.Item(25, 7).Value = 100101
$ws.Cells.Item(25, 8).Value = "Berries"
$ws.Cells.Item(25, 9).Value = 100101001
$ws.Cells.Item(25, 10).Value = "Arándano (blue)"
$ws.Cells.Item(25, 11).Value = "Sin especificar"
$ws.Cells.Item(25, 12).Value = "Primera"
$ws.Cells.Item(25, 13).Value = 150
$ws.Cells.Item(25, 14).Value = 3500
$ws.Cells.Item(25, 15).Value = 3500
$ws.Cells.Item(25, 16).Value = 3500
$ws.Cells.Item(25, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(25, 18).Value = "Provincia de Linares"
$ws.Cells.Item(25, 19).Value = 1750
$ws.Cells.Item(25, 20).Value = 2
